$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.047.85"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.06%  "

$ws.Range("D3").Value = "'2.057.56"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.42%  "

$ws.Range("D5").Value = "'246.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.93%  "

$ws.Range("D6").Value = "'0.660"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.28%  "

$ws.Range("D7").Value = "'58.81"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.64%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("E9").Value = "  -2.36%  "

$ws.Range("D10").Value = "'0.0784"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.88%  "

$ws.Range("E11").Value = "  +2.46%  "

$ws.Range("D12").Value = "'15.48"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.00%  "

$ws.Range("D13").Value = "'0.890"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.90%  "

$ws.Range("D14").Value = "'2.357.72"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.45%  "

$ws.Range("D15").Value = "'5.73"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.28%  "

$ws.Range("D16").Value = "'2.068.54"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.13%  "

$ws.Range("D17").Value = "'18.26"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.16%  "

$ws.Range("D18").Value = "'37.024.79"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.03%  "

$ws.Range("D19").Value = "'73.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.87%  "

$ws.Range("D20").Value = "'0.0₃0897"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.56%  "

$ws.Range("D21").Value = "'5.46"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.86%  "

$ws.Range("D22").Value = "'238.59"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.33%  "

$ws.Range("E23").Value = "  -0.02%  "

$ws.Range("E24").Value = "  +1.62%  "

$ws.Range("D25").Value = "'10.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +8.18%  "

$ws.Range("D26").Value = "'170.59"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.13%  "

$ws.Range("E27").Value = "  -2.01%  "

$ws.Range("D28").Value = "'20.16"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.40%  "

$ws.Range("D29").Value = "'5.53"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +15.48%  "

$ws.Range("E30").Value = "  -0.93%  "

$ws.Range("D31").Value = "'1.12"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.43%  "

$ws.Range("D32").Value = "'4.71"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.50%  "

$ws.Range("D33").Value = "'0.0620"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.82%  "

$ws.Range("D34").Value = "'2.36"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.45%  "

$ws.Range("E35").Value = "  -0.12%  "

$ws.Range("D36").Value = "'1.84"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.62%  "

$ws.Range("D37").Value = "'0.0847"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.75%  "

$ws.Range("D38").Value = "'1.34"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.04%  "

$ws.Range("D39").Value = "'5.25"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.25%  "

$ws.Range("E40").Value = "  -1.32%  "

$ws.Range("E41").Value = "  +0.70%  "

$ws.Range("D42").Value = "'1.17"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.38%  "

$ws.Range("D43").Value = "'0.0965"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -10.03%  "

$ws.Range("D44").Value = "'97.88"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.02%  "

$ws.Range("D45").Value = "'17.04"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.57%  "

$ws.Range("D46").Value = "'1.303.92"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.19%  "

$ws.Range("D47").Value = "'2.38"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.73%  "

$ws.Range("E48").Value = "  -0.43%  "

$ws.Range("E49").Value = "  +0.30%  "

$ws.Range("D50").Value = "'2.247.34"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.82%  "

$ws.Range("D51").Value = "'44.78"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.36%  "
